# Refresh the GSC "Coverage" export: the crawl window rolled forward by
# one day, so the oldest date row (2025-10-15 -- the first data row under
# the header on the "Chart" sheet) drops off the front of the report and
# every later row shifts up to take its place. This is the "updated main
# GSC export data" refresh described by the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date in the series; deleting the whole row
# shifts rows 3:89 up into 2:88 (real cell values, not just formatting),
# matching how the upstream exporter re-wrote the sheet.
$ws.Range("A2").EntireRow.Delete()
